# Upgrade selenium version - add new login credentials row to the
# "Login Credentials" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login Credentials")

# Find the next empty row after the existing data (row 42 -> new row 43)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$ws.Cells.Item($newRow, 1).Value = "fittd@gmail.com"
$ws.Cells.Item($newRow, 2).Value = "SUu6bq349$"
$ws.Cells.Item($newRow, 3).Value = "Birjesh"
$ws.Cells.Item($newRow, 4).Value = "Bandopadhyay D"
$ws.Cells.Item($newRow, 5).Value = "Birjesh Bandopadhyay D"
$ws.Cells.Item($newRow, 6).Value = ","

$wb.Save()
